$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 5 - España
$ws.Range("B5").Value = 236199
$ws.Range("C5").Value = 9570
$ws.Range("D5").Value = 127609
$ws.Range("E5").Value = 85069
$ws.Range("G5").Value = 331
$ws.Range("H5").Value = 23521

# Row 8 - Alemania
$ws.Range("B8").Value = 157781
$ws.Range("C8").Value = 11
$ws.Range("E8").Value = 37305

# Row 11 - Iran
$ws.Range("B11").Value = 91472
$ws.Range("C11").Value = 991
$ws.Range("D11").Value = 70933
$ws.Range("E11").Value = 14733
$ws.Range("F11").Value = 3011
$ws.Range("G11").Value = 96
$ws.Range("H11").Value = 5806

# Row 18 - Suiza
$ws.Range("E18").Value = 5621
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 1640
